$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.325823666666667
$ws.Range("H2").Value = 18.977471
$ws.Range("I2").Value = 0.03655606111596783
$ws.Range("J2").Value = 0.0371094793300862
$ws.Range("M2").Value = 6.038588
$ws.Range("N2").Value = 18.115764
$ws.Range("O2").Value = 0.4832124175152646
$ws.Range("P2").Value = 0.5250116720691949
$ws.Range("Q2").Value = 38.19904288364933
$ws.Range("R2").Value = 343.791385952844
$ws.Range("S2").Value = 0.01766434266668258
$ws.Range("T2").Value = 0.01948290979270578
$ws.Range("G3").Value = 6.325823666666667
$ws.Range("H3").Value = 18.977471
$ws.Range("I3").Value = 0.03655606111596783
$ws.Range("J3").Value = 0.0371094793300862
$ws.Range("O3").Value = 0.157682833439183
$ws.Range("P3").Value = 0.1713228489992161
$ws.Range("Q3").Value = 12.46518735493467
$ws.Range("R3").Value = 112.186686194412
$ws.Range("S3").Value = 0.00576426329614175
$ws.Range("T3").Value = 0.006357701723707888
$ws.Range("G4").Value = 6.325823666666667
$ws.Range("H4").Value = 18.977471
$ws.Range("I4").Value = 0.03655606111596783
$ws.Range("J4").Value = 0.0371094793300862
$ws.Range("M4").Value = 0.6106236666666667
$ws.Range("N4").Value = 1.831871
$ws.Range("O4").Value = 0.04886257154189607
$ws.Range("P4").Value = 0.05308932357062435
$ws.Range("Q4").Value = 3.862697642026778
$ws.Range("R4").Value = 34.764278778241
$ws.Range("S4").Value = 0.001786223151568903
$ws.Range("T4").Value = 0.001970117155692343
$ws.Range("G5").Value = 6.325823666666667
$ws.Range("H5").Value = 18.977471
$ws.Range("I5").Value = 0.03655606111596783
$ws.Range("J5").Value = 0.0371094793300862
$ws.Range("M5").Value = 2.98482
$ws.Range("N5").Value = 5.96964
$ws.Range("O5").Value = 0.2388475729836035
$ws.Range("P5").Value = 0.173005713590172
$ws.Range("Q5").Value = 18.88144499674
$ws.Range("R5").Value = 113.28866998044
$ws.Range("S5").Value = 0.008731326475389194
$ws.Range("T5").Value = 0.006420151952461301
$ws.Range("G6").Value = 6.325823666666667
$ws.Range("H6").Value = 18.977471
$ws.Range("I6").Value = 0.03655606111596783
$ws.Range("J6").Value = 0.0371094793300862
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.892201
$ws.Range("N6").Value = 2.676603
$ws.Range("O6").Value = 0.07139460452005281
$ws.Range("P6").Value = 0.07757044177079274
$ws.Range("Q6").Value = 5.643906201223667
$ws.Range("R6").Value = 50.795155811013
$ws.Range("S6").Value = 0.002609905526185403
$ws.Range("T6").Value = 0.002878598705518889
$ws.Range("I7").Value = 0.4487538134336191
$ws.Range("J7").Value = 0.4555474483720607
$ws.Range("M7").Value = 6.038588
$ws.Range("N7").Value = 18.115764
$ws.Range("O7").Value = 0.4832124175152646
$ws.Range("P7").Value = 0.5250116720691949
$ws.Range("Q7").Value = 468.922680405092
$ws.Range("R7").Value = 4220.304123645828
$ws.Range("S7").Value = 0.2168434150584531
$ws.Range("T7").Value = 0.2391677275766708
$ws.Range("I8").Value = 0.4487538134336191
$ws.Range("J8").Value = 0.4555474483720607
$ws.Range("O8").Value = 0.157682833439183
$ws.Range("P8").Value = 0.1713228489992161
$ws.Range("S8").Value = 0.07076077281885158
$ws.Range("T8").Value = 0.07804568670942473
$ws.Range("I9").Value = 0.4487538134336191
$ws.Range("J9").Value = 0.4555474483720607
$ws.Range("M9").Value = 0.6106236666666667
$ws.Range("N9").Value = 1.831871
$ws.Range("O9").Value = 0.04886257154189607
$ws.Range("P9").Value = 0.05308932357062435
$ws.Range("Q9").Value = 47.41758942522966
$ws.Range("R9").Value = 426.758304827067
$ws.Range("S9").Value = 0.0219272653135989
$ws.Range("T9").Value = 0.02418470588839662
$ws.Range("I10").Value = 0.4487538134336191
$ws.Range("J10").Value = 0.4555474483720607
$ws.Range("M10").Value = 2.98482
$ws.Range("N10").Value = 5.96964
$ws.Range("O10").Value = 0.2388475729836035
$ws.Range("P10").Value = 0.173005713590172
$ws.Range("Q10").Value = 231.78428383038
$ws.Range("R10").Value = 1390.70570298228
$ws.Range("S10").Value = 0.1071837592057567
$ws.Range("T10").Value = 0.07881231137979039
$ws.Range("I11").Value = 0.4487538134336191
$ws.Range("J11").Value = 0.4555474483720607
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.892201
$ws.Range("N11").Value = 2.676603
$ws.Range("O11").Value = 0.07139460452005281
$ws.Range("P11").Value = 0.07757044177079274
$ws.Range("Q11").Value = 69.283296754159
$ws.Range("R11").Value = 623.549670787431
$ws.Range("S11").Value = 0.0320386010369588
$ws.Range("T11").Value = 0.03533701681777814
$ws.Range("G12").Value = 31.09589533333333
$ws.Range("H12").Value = 93.28768600000001
$ws.Range("I12").Value = 0.1796988835226367
$ws.Range("J12").Value = 0.1824193252814652
$ws.Range("M12").Value = 6.038588
$ws.Range("N12").Value = 18.115764
$ws.Range("O12").Value = 0.4832124175152646
$ws.Range("P12").Value = 0.5250116720691949
$ws.Range("Q12").Value = 187.7753004091227
$ws.Range("R12").Value = 1689.977703682104
$ws.Range("S12").Value = 0.08683273193176724
$ws.Range("T12").Value = 0.09577227498375637
$ws.Range("G13").Value = 31.09589533333333
$ws.Range("H13").Value = 93.28768600000001
$ws.Range("I13").Value = 0.1796988835226367
$ws.Range("J13").Value = 0.1824193252814652
$ws.Range("O13").Value = 0.157682833439183
$ws.Range("P13").Value = 0.1713228489992161
$ws.Range("Q13").Value = 61.27520805582134
$ws.Range("R13").Value = 551.476872502392
$ws.Range("S13").Value = 0.02833542911970708
$ws.Range("T13").Value = 0.03125259851973533
$ws.Range("G14").Value = 31.09589533333333
$ws.Range("H14").Value = 93.28768600000001
$ws.Range("I14").Value = 0.1796988835226367
$ws.Range("J14").Value = 0.1824193252814652
$ws.Range("M14").Value = 0.6106236666666667
$ws.Range("N14").Value = 1.831871
$ws.Range("O14").Value = 0.04886257154189607
$ws.Range("P14").Value = 0.05308932357062435
$ws.Range("Q14").Value = 18.98788962672289
$ws.Range("R14").Value = 170.891006640506
$ws.Range("S14").Value = 0.008780549552123686
$ws.Range("T14").Value = 0.009684518585402679
$ws.Range("G15").Value = 31.09589533333333
$ws.Range("H15").Value = 93.28768600000001
$ws.Range("I15").Value = 0.1796988835226367
$ws.Range("J15").Value = 0.1824193252814652
$ws.Range("M15").Value = 2.98482
$ws.Range("N15").Value = 5.96964
$ws.Range("O15").Value = 0.2388475729836035
$ws.Range("P15").Value = 0.173005713590172
$ws.Range("Q15").Value = 92.81565030884001
$ws.Range("R15").Value = 556.89390185304
$ws.Range("S15").Value = 0.04292064219724503
$ws.Range("T15").Value = 0.03155958554295758
$ws.Range("G16").Value = 31.09589533333333
$ws.Range("H16").Value = 93.28768600000001
$ws.Range("I16").Value = 0.1796988835226367
$ws.Range("J16").Value = 0.1824193252814652
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.892201
$ws.Range("N16").Value = 2.676603
$ws.Range("O16").Value = 0.07139460452005281
$ws.Range("P16").Value = 0.07757044177079274
$ws.Range("Q16").Value = 27.74378891229534
$ws.Range("R16").Value = 249.694100210658
$ws.Range("S16").Value = 0.01282953072179368
$ws.Range("T16").Value = 0.01415034764961319
$ws.Range("G17").Value = 7.741899500000001
$ws.Range("H17").Value = 15.483799
$ws.Range("I17").Value = 0.04473936773909665
$ws.Range("J17").Value = 0.0302777814252336
$ws.Range("M17").Value = 6.038588
$ws.Range("N17").Value = 18.115764
$ws.Range("O17").Value = 0.4832124175152646
$ws.Range("P17").Value = 0.5250116720691949
$ws.Range("Q17").Value = 46.750141417906
$ws.Range("R17").Value = 280.500848507436
$ws.Range("S17").Value = 0.02161861804331333
$ws.Range("T17").Value = 0.0158961886526075
$ws.Range("G18").Value = 7.741899500000001
$ws.Range("H18").Value = 15.483799
$ws.Range("I18").Value = 0.04473936773909665
$ws.Range("J18").Value = 0.0302777814252336
$ws.Range("O18").Value = 0.157682833439183
$ws.Range("P18").Value = 0.1713228489992161
$ws.Range("Q18").Value = 15.255598770338
$ws.Range("R18").Value = 91.533592622028
$ws.Range("S18").Value = 0.007054630271378335
$ws.Range("T18").Value = 0.005187275775146566
$ws.Range("G19").Value = 7.741899500000001
$ws.Range("H19").Value = 15.483799
$ws.Range("I19").Value = 0.04473936773909665
$ws.Range("J19").Value = 0.0302777814252336
$ws.Range("M19").Value = 0.6106236666666667
$ws.Range("N19").Value = 1.831871
$ws.Range("O19").Value = 0.04886257154189607
$ws.Range("P19").Value = 0.05308932357062435
$ws.Range("Q19").Value = 4.727387059654834
$ws.Range("R19").Value = 28.364322357929
$ws.Range("S19").Value = 0.002186080556890807
$ws.Range("T19").Value = 0.001607426935084866
$ws.Range("G20").Value = 7.741899500000001
$ws.Range("H20").Value = 15.483799
$ws.Range("I20").Value = 0.04473936773909665
$ws.Range("J20").Value = 0.0302777814252336
$ws.Range("M20").Value = 2.98482
$ws.Range("N20").Value = 5.96964
$ws.Range("O20").Value = 0.2388475729836035
$ws.Range("P20").Value = 0.173005713590172
$ws.Range("Q20").Value = 23.10817646559
$ws.Range("R20").Value = 92.43270586236001
$ws.Range("S20").Value = 0.01068588940130416
$ws.Range("T20").Value = 0.005238229181399794
$ws.Range("G21").Value = 7.741899500000001
$ws.Range("H21").Value = 15.483799
$ws.Range("I21").Value = 0.04473936773909665
$ws.Range("J21").Value = 0.0302777814252336
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.892201
$ws.Range("N21").Value = 2.676603
$ws.Range("O21").Value = 0.07139460452005281
$ws.Range("P21").Value = 0.07757044177079274
$ws.Range("Q21").Value = 6.907330475799501
$ws.Range("R21").Value = 41.44398285479701
$ws.Range("S21").Value = 0.003194149466210014
$ws.Range("T21").Value = 0.002348660880994873
$ws.Range("G22").Value = 50.22647733333333
$ws.Range("H22").Value = 150.679432
$ws.Range("I22").Value = 0.2902518741886796
$ws.Range("J22").Value = 0.2946459655911543
$ws.Range("M22").Value = 6.038588
$ws.Range("N22").Value = 18.115764
$ws.Range("O22").Value = 0.4832124175152646
$ws.Range("P22").Value = 0.5250116720691949
$ws.Range("Q22").Value = 303.2970033073386
$ws.Range("R22").Value = 2729.673029766047
$ws.Range("S22").Value = 0.1402533098150483
$ws.Range("T22").Value = 0.1546925710634544
$ws.Range("G23").Value = 50.22647733333333
$ws.Range("H23").Value = 150.679432
$ws.Range("I23").Value = 0.2902518741886796
$ws.Range("J23").Value = 0.2946459655911543
$ws.Range("O23").Value = 0.157682833439183
$ws.Range("P23").Value = 0.1713228489992161
$ws.Range("Q23").Value = 98.97247902078932
$ws.Range("R23").Value = 890.7523111871039
$ws.Range("S23").Value = 0.04576773793310428
$ws.Range("T23").Value = 0.05047958627120154
$ws.Range("G24").Value = 50.22647733333333
$ws.Range("H24").Value = 150.679432
$ws.Range("I24").Value = 0.2902518741886796
$ws.Range("J24").Value = 0.2946459655911543
$ws.Range("M24").Value = 0.6106236666666667
$ws.Range("N24").Value = 1.831871
$ws.Range("O24").Value = 0.04886257154189607
$ws.Range("P24").Value = 0.05308932357062435
$ws.Range("Q24").Value = 30.66947575303022
$ws.Range("R24").Value = 276.025281777272
$ws.Range("S24").Value = 0.01418245296771378
$ws.Range("T24").Value = 0.01564255500604784
$ws.Range("G25").Value = 50.22647733333333
$ws.Range("H25").Value = 150.679432
$ws.Range("I25").Value = 0.2902518741886796
$ws.Range("J25").Value = 0.2946459655911543
$ws.Range("M25").Value = 2.98482
$ws.Range("N25").Value = 5.96964
$ws.Range("O25").Value = 0.2388475729836035
$ws.Range("P25").Value = 0.173005713590172
$ws.Range("Q25").Value = 149.91699407408
$ws.Range("R25").Value = 899.50196444448
$ws.Range("S25").Value = 0.06932595570390836
$ws.Range("T25").Value = 0.05097543553356292
$ws.Range("G26").Value = 50.22647733333333
$ws.Range("H26").Value = 150.679432
$ws.Range("I26").Value = 0.2902518741886796
$ws.Range("J26").Value = 0.2946459655911543
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 0.892201
$ws.Range("N26").Value = 2.676603
$ws.Range("O26").Value = 0.07139460452005281
$ws.Range("P26").Value = 0.07757044177079274
$ws.Range("Q26").Value = 44.81211330327733
$ws.Range("R26").Value = 403.309019729496
$ws.Range("S26").Value = 0.0207224177689049
$ws.Range("T26").Value = 0.02285581771688764
